# Updates cryptos list values (price & volume-1h columns, plus two row swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.734.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.103.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.05%  "

$ws.Range("E4").Value = "  +0.51%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "346.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.27%  "

$ws.Range("E6").Value = "  +0.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5195"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.52%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4435"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.22"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09353"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.176"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.157.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.24%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.819"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.34%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.299"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001160"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.81%  "

$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06674"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.310"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.10%  "

$ws.Range("E22").Value = "  +0.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.791.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.323"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.392.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.543"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.11%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.145"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.72%  "

$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.792"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.30%  "

$ws.Range("E33").Value = "  -1.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.224"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.29%  "

$ws.Range("E35").Value = "  -0.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.333"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02590"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06774"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.7024"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.335"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2233"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6829"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.90%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.359"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.80%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.637"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000355"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.222"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.221"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.50%  "
